$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-10 down to 8-11.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly price entry.
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 45005
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100114002
$ws.Cells.Item(7, 7).Value = "Camote"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 200
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 972
$ws.Cells.Item(7, 17).Value = 18
$ws.Cells.Item(7, 18).Value = "Hortaliza"
